# Re-run SGNN to annotate dialog acts following clean up work to the original transcripts.
# Update DAMSLTag (column I) and DialogAct (column J) values for the affected rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(13, 9).Value = 'sd'
$ws.Cells.Item(13, 10).Value = 'Statement-non-opinion'
$ws.Cells.Item(44, 9).Value = 'sd'
$ws.Cells.Item(44, 10).Value = 'Statement-non-opinion'
$ws.Cells.Item(46, 9).Value = 'sv'
$ws.Cells.Item(46, 10).Value = 'Statement-opinion'
$ws.Cells.Item(64, 9).Value = 'sd'
$ws.Cells.Item(64, 10).Value = 'Statement-non-opinion'
$ws.Cells.Item(86, 9).Value = 'sv'
$ws.Cells.Item(86, 10).Value = 'Statement-opinion'
$ws.Cells.Item(99, 9).Value = 'ba'
$ws.Cells.Item(99, 10).Value = 'Appreciation'
$ws.Cells.Item(105, 9).Value = 'ba'
$ws.Cells.Item(105, 10).Value = 'Appreciation'
$ws.Cells.Item(111, 9).Value = 'aa'
$ws.Cells.Item(111, 10).Value = 'Agree/Accept'
$ws.Cells.Item(117, 9).Value = 'ba'
$ws.Cells.Item(117, 10).Value = 'Appreciation'
$ws.Cells.Item(130, 9).Value = 'b'
$ws.Cells.Item(130, 10).Value = 'Acknowledge (Backchannel)'
$ws.Cells.Item(135, 9).Value = 'qy'
$ws.Cells.Item(135, 10).Value = 'Yes-No-Question'
$ws.Cells.Item(136, 9).Value = 'sv'
$ws.Cells.Item(136, 10).Value = 'Statement-opinion'
$ws.Cells.Item(139, 9).Value = 'sv'
$ws.Cells.Item(139, 10).Value = 'Statement-opinion'
$ws.Cells.Item(147, 9).Value = 'sv'
$ws.Cells.Item(147, 10).Value = 'Statement-opinion'
$ws.Cells.Item(148, 9).Value = 'sv'
$ws.Cells.Item(148, 10).Value = 'Statement-opinion'
$ws.Cells.Item(149, 9).Value = 'qy'
$ws.Cells.Item(149, 10).Value = 'Yes-No-Question'
$ws.Cells.Item(159, 9).Value = 'sv'
$ws.Cells.Item(159, 10).Value = 'Statement-opinion'
$ws.Cells.Item(162, 9).Value = 'sd'
$ws.Cells.Item(162, 10).Value = 'Statement-non-opinion'
$ws.Cells.Item(169, 9).Value = '%'
$ws.Cells.Item(169, 10).Value = 'Uninterpretable'
$ws.Cells.Item(175, 9).Value = 'sd'
$ws.Cells.Item(175, 10).Value = 'Statement-non-opinion'
$ws.Cells.Item(179, 9).Value = 'sv'
$ws.Cells.Item(179, 10).Value = 'Statement-opinion'
$ws.Cells.Item(201, 9).Value = 'b'
$ws.Cells.Item(201, 10).Value = 'Acknowledge (Backchannel)'
$ws.Cells.Item(208, 9).Value = 'qy'
$ws.Cells.Item(208, 10).Value = 'Yes-No-Question'
$ws.Cells.Item(209, 9).Value = '%'
$ws.Cells.Item(209, 10).Value = 'Uninterpretable'
$ws.Cells.Item(211, 9).Value = 'ba'
$ws.Cells.Item(211, 10).Value = 'Appreciation'
$ws.Cells.Item(216, 9).Value = 'aa'
$ws.Cells.Item(216, 10).Value = 'Agree/Accept'
$ws.Cells.Item(223, 9).Value = 'sv'
$ws.Cells.Item(223, 10).Value = 'Statement-opinion'
$ws.Cells.Item(234, 9).Value = 'aa'
$ws.Cells.Item(234, 10).Value = 'Agree/Accept'
$ws.Cells.Item(247, 9).Value = 'aa'
$ws.Cells.Item(247, 10).Value = 'Agree/Accept'
$ws.Cells.Item(250, 9).Value = 'sd'
$ws.Cells.Item(250, 10).Value = 'Statement-non-opinion'
$ws.Cells.Item(252, 9).Value = 'sd'
$ws.Cells.Item(252, 10).Value = 'Statement-non-opinion'
$ws.Cells.Item(253, 9).Value = 'sd'
$ws.Cells.Item(253, 10).Value = 'Statement-non-opinion'
$ws.Cells.Item(257, 9).Value = 'aa'
$ws.Cells.Item(257, 10).Value = 'Agree/Accept'
$ws.Cells.Item(260, 9).Value = 'ba'
$ws.Cells.Item(260, 10).Value = 'Appreciation'
$ws.Cells.Item(261, 9).Value = 'b'
$ws.Cells.Item(261, 10).Value = 'Acknowledge (Backchannel)'
$ws.Cells.Item(262, 9).Value = 'ba'
$ws.Cells.Item(262, 10).Value = 'Appreciation'
$ws.Cells.Item(268, 9).Value = 'b'
$ws.Cells.Item(268, 10).Value = 'Acknowledge (Backchannel)'
$ws.Cells.Item(284, 9).Value = 'sd'
$ws.Cells.Item(284, 10).Value = 'Statement-non-opinion'
$ws.Cells.Item(287, 9).Value = 'sd'
$ws.Cells.Item(287, 10).Value = 'Statement-non-opinion'
$ws.Cells.Item(306, 9).Value = '%'
$ws.Cells.Item(306, 10).Value = 'Uninterpretable'
$ws.Cells.Item(308, 9).Value = 'sv'
$ws.Cells.Item(308, 10).Value = 'Statement-opinion'
$ws.Cells.Item(320, 9).Value = 'ba'
$ws.Cells.Item(320, 10).Value = 'Appreciation'
$ws.Cells.Item(337, 9).Value = 'sv'
$ws.Cells.Item(337, 10).Value = 'Statement-opinion'
$ws.Cells.Item(338, 9).Value = 'sd'
$ws.Cells.Item(338, 10).Value = 'Statement-non-opinion'
$ws.Cells.Item(341, 9).Value = 'sd'
$ws.Cells.Item(341, 10).Value = 'Statement-non-opinion'
$ws.Cells.Item(351, 9).Value = 'ba'
$ws.Cells.Item(351, 10).Value = 'Appreciation'
$ws.Cells.Item(352, 9).Value = 'sd'
$ws.Cells.Item(352, 10).Value = 'Statement-non-opinion'
$ws.Cells.Item(355, 9).Value = 'aa'
$ws.Cells.Item(355, 10).Value = 'Agree/Accept'
$ws.Cells.Item(359, 9).Value = 'sd'
$ws.Cells.Item(359, 10).Value = 'Statement-non-opinion'
$ws.Cells.Item(362, 9).Value = 'sv'
$ws.Cells.Item(362, 10).Value = 'Statement-opinion'
$ws.Cells.Item(363, 9).Value = 'aa'
$ws.Cells.Item(363, 10).Value = 'Agree/Accept'
$ws.Cells.Item(370, 9).Value = 'sv'
$ws.Cells.Item(370, 10).Value = 'Statement-opinion'
$ws.Cells.Item(389, 9).Value = 'b'
$ws.Cells.Item(389, 10).Value = 'Acknowledge (Backchannel)'
$ws.Cells.Item(393, 9).Value = 'sd'
$ws.Cells.Item(393, 10).Value = 'Statement-non-opinion'
$ws.Cells.Item(394, 9).Value = 'sd'
$ws.Cells.Item(394, 10).Value = 'Statement-non-opinion'
$ws.Cells.Item(405, 9).Value = 'sd'
$ws.Cells.Item(405, 10).Value = 'Statement-non-opinion'
$ws.Cells.Item(407, 9).Value = 'sd'
$ws.Cells.Item(407, 10).Value = 'Statement-non-opinion'
$ws.Cells.Item(412, 9).Value = 'sv'
$ws.Cells.Item(412, 10).Value = 'Statement-opinion'
$ws.Cells.Item(417, 9).Value = 'sd'
$ws.Cells.Item(417, 10).Value = 'Statement-non-opinion'
$ws.Cells.Item(440, 9).Value = 'sd'
$ws.Cells.Item(440, 10).Value = 'Statement-non-opinion'
$ws.Cells.Item(444, 9).Value = 'sv'
$ws.Cells.Item(444, 10).Value = 'Statement-opinion'
$ws.Cells.Item(446, 9).Value = 'sd'
$ws.Cells.Item(446, 10).Value = 'Statement-non-opinion'
$ws.Cells.Item(449, 9).Value = '%'
$ws.Cells.Item(449, 10).Value = 'Uninterpretable'
$ws.Cells.Item(454, 9).Value = 'ba'
$ws.Cells.Item(454, 10).Value = 'Appreciation'
$ws.Cells.Item(456, 9).Value = 'sv'
$ws.Cells.Item(456, 10).Value = 'Statement-opinion'
$ws.Cells.Item(459, 9).Value = 'sv'
$ws.Cells.Item(459, 10).Value = 'Statement-opinion'
$ws.Cells.Item(460, 9).Value = 'sv'
$ws.Cells.Item(460, 10).Value = 'Statement-opinion'
$ws.Cells.Item(466, 9).Value = 'sv'
$ws.Cells.Item(466, 10).Value = 'Statement-opinion'
$ws.Cells.Item(477, 9).Value = 'b'
$ws.Cells.Item(477, 10).Value = 'Acknowledge (Backchannel)'
$ws.Cells.Item(485, 9).Value = 'sv'
$ws.Cells.Item(485, 10).Value = 'Statement-opinion'
$ws.Cells.Item(486, 9).Value = 'sd'
$ws.Cells.Item(486, 10).Value = 'Statement-non-opinion'
$ws.Cells.Item(499, 9).Value = 'aa'
$ws.Cells.Item(499, 10).Value = 'Agree/Accept'
$ws.Cells.Item(507, 9).Value = 'sd'
$ws.Cells.Item(507, 10).Value = 'Statement-non-opinion'
$ws.Cells.Item(509, 9).Value = 'sd'
$ws.Cells.Item(509, 10).Value = 'Statement-non-opinion'
$ws.Cells.Item(512, 9).Value = 'aa'
$ws.Cells.Item(512, 10).Value = 'Agree/Accept'
$ws.Cells.Item(516, 9).Value = '%'
$ws.Cells.Item(516, 10).Value = 'Uninterpretable'
$ws.Cells.Item(523, 9).Value = 'sd'
$ws.Cells.Item(523, 10).Value = 'Statement-non-opinion'
$ws.Cells.Item(524, 9).Value = 'b'
$ws.Cells.Item(524, 10).Value = 'Acknowledge (Backchannel)'
$ws.Cells.Item(531, 9).Value = 'b'
$ws.Cells.Item(531, 10).Value = 'Acknowledge (Backchannel)'
$ws.Cells.Item(535, 9).Value = 'b'
$ws.Cells.Item(535, 10).Value = 'Acknowledge (Backchannel)'
$ws.Cells.Item(538, 9).Value = 'aa'
$ws.Cells.Item(538, 10).Value = 'Agree/Accept'
$ws.Cells.Item(542, 9).Value = 'aa'
$ws.Cells.Item(542, 10).Value = 'Agree/Accept'
$ws.Cells.Item(545, 9).Value = 'sv'
$ws.Cells.Item(545, 10).Value = 'Statement-opinion'
$ws.Cells.Item(549, 9).Value = 'aa'
$ws.Cells.Item(549, 10).Value = 'Agree/Accept'
$ws.Cells.Item(555, 9).Value = 'b'
$ws.Cells.Item(555, 10).Value = 'Acknowledge (Backchannel)'
$ws.Cells.Item(556, 9).Value = 'sd'
$ws.Cells.Item(556, 10).Value = 'Statement-non-opinion'
$ws.Cells.Item(564, 9).Value = 'sd'
$ws.Cells.Item(564, 10).Value = 'Statement-non-opinion'
$ws.Cells.Item(572, 9).Value = 'aa'
$ws.Cells.Item(572, 10).Value = 'Agree/Accept'
